$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("Overview")
$ws1.Range("A2").Value = "80ceb9b9-d06f-4040-b9f6-55dbbbfa287d.md"
$ws1.Range("A3").Value = "f397695c-0ec3-4d6a-a189-1b31f3b14007.md"
$ws1.Range("B3").Value = "Handed back: in sync with en-US"
$ws1.Range("C3").Value = "Handed back: in sync with en-US"

$ws2 = $wb.Worksheets.Item("zh-cn")
$ws2.Range("A2").Value = "80ceb9b9-d06f-4040-b9f6-55dbbbfa287d.md"
$ws2.Range("C2").Value = "80ceb9b9-d06f-4040-b9f6-55dbbbfa287d.24877f15d3eefb0c26507cc2860096079cab7782.zh-cn.xlf"
$ws2.Range("E2").Value = "80ceb9b9-d06f-4040-b9f6-55dbbbfa287d.md"
$ws2.Range("F2").Value = "80ceb9b9-d06f-4040-b9f6-55dbbbfa287d.24877f15d3eefb0c26507cc2860096079cab7782.zh-cn.xlf"
$ws2.Range("G2").Value = "2016-03-04 06:44:56"
$ws2.Range("A3").Value = "f397695c-0ec3-4d6a-a189-1b31f3b14007.md"
$ws2.Range("B3").Value = "Handed back: in sync with en-US"
$ws2.Range("C3").Value = "f397695c-0ec3-4d6a-a189-1b31f3b14007.d03e07658630811d43ed625d04768daa097b053c.zh-cn.xlf"
$ws2.Range("E3").Value = "f397695c-0ec3-4d6a-a189-1b31f3b14007.md"
$ws2.Range("F3").Value = "f397695c-0ec3-4d6a-a189-1b31f3b14007.d03e07658630811d43ed625d04768daa097b053c.zh-cn.xlf"
$ws2.Range("G3").Value = "2016-03-04 06:44:56"

$ws3 = $wb.Worksheets.Item("de-de")
$ws3.Range("A2").Value = "80ceb9b9-d06f-4040-b9f6-55dbbbfa287d.md"
$ws3.Range("C2").Value = "80ceb9b9-d06f-4040-b9f6-55dbbbfa287d.24877f15d3eefb0c26507cc2860096079cab7782.de-de.xlf"
$ws3.Range("E2").Value = "80ceb9b9-d06f-4040-b9f6-55dbbbfa287d.md"
$ws3.Range("F2").Value = "80ceb9b9-d06f-4040-b9f6-55dbbbfa287d.24877f15d3eefb0c26507cc2860096079cab7782.de-de.xlf"
$ws3.Range("G2").Value = "2016-03-04 06:45:26"
$ws3.Range("A3").Value = "f397695c-0ec3-4d6a-a189-1b31f3b14007.md"
$ws3.Range("B3").Value = "Handed back: in sync with en-US"
$ws3.Range("C3").Value = "f397695c-0ec3-4d6a-a189-1b31f3b14007.d03e07658630811d43ed625d04768daa097b053c.de-de.xlf"
$ws3.Range("E3").Value = "f397695c-0ec3-4d6a-a189-1b31f3b14007.md"
$ws3.Range("F3").Value = "f397695c-0ec3-4d6a-a189-1b31f3b14007.d03e07658630811d43ed625d04768daa097b053c.de-de.xlf"
$ws3.Range("G3").Value = "2016-03-04 06:45:26"
